# Auto-generated Excel COM-interop script
# Updates scheduled-runner price/profit snapshots across all 8 Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per Sheets/Gilgamesh_Profits.xlsx diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 499.83334
$ws.Range("I33").Value = 499.83334
$ws.Range("K33").Value = 499.83334
$ws.Range("M33").Value = -270.83334
$ws.Range("H40").Value = 6111.1113
$ws.Range("I40").Value = 5250
$ws.Range("K40").Value = 5250
$ws.Range("M40").Value = -5075
$ws.Range("H74").Value = 12712.863
$ws.Range("I74").Value = 13059.15
$ws.Range("K74").Value = 13059.15
$ws.Range("M74").Value = -12123.15
$ws.Range("H76").Value = 3999.2
$ws.Range("I76").Value = 3999
$ws.Range("K76").Value = 3999
$ws.Range("M76").Value = -3684
$ws.Range("H77").Value = 12712.863
$ws.Range("I77").Value = 13059.15
$ws.Range("K77").Value = 65295.75
$ws.Range("M77").Value = -60615.75
$ws.Range("H79").Value = 3999.2
$ws.Range("I79").Value = 3999
$ws.Range("K79").Value = 3999
$ws.Range("M79").Value = -2907
$ws.Range("H88").Value = 5264843
$ws.Range("J88").Value = 2009.3077
$ws.Range("L88").Value = 2009.3077
$ws.Range("N88").Value = -2821.3077
$ws.Range("H91").Value = 5264843
$ws.Range("J91").Value = 2009.3077
$ws.Range("L91").Value = 2009.3077
$ws.Range("N91").Value = -4817.3077
$ws.Range("H106").Value = 1995.3334
$ws.Range("I106").Value = 1995.3334
$ws.Range("K106").Value = 1995.3334
$ws.Range("M106").Value = -1364.3334
$ws.Range("H107").Value = 1219.8667
$ws.Range("J107").Value = 1678.8
$ws.Range("L107").Value = 1678.8
$ws.Range("N107").Value = -5518.8
$ws.Range("H137").Value = 6252725
$ws.Range("I137").Value = 6252725
$ws.Range("K137").Value = 18758175
$ws.Range("M137").Value = -18755625
$ws.Range("H138").Value = 5702.5757
$ws.Range("I138").Value = 4249.8335
$ws.Range("J138").Value = 6025.407
$ws.Range("K138").Value = 12749.5005
$ws.Range("L138").Value = 18076.221
$ws.Range("M138").Value = -7609.500499999998
$ws.Range("N138").Value = -28356.221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1842173
$ws.Range("I32").Value = 824909.1
$ws.Range("K32").Value = 824909.1
$ws.Range("M32").Value = -824622.1
$ws.Range("H132").Value = 1548.84
$ws.Range("I132").Value = 1177.1904
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 3531.5712
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -1001.5712
$ws.Range("N132").Value = -15560
$ws.Range("H135").Value = 94400
$ws.Range("J135").Value = 94400
$ws.Range("L135").Value = 94400
$ws.Range("N135").Value = -104540

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 50
$ws.Range("I8").Value = 50
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 50
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 90
$ws.Range("N8").ClearContents()
$ws.Range("H86").Value = 2841.5173
$ws.Range("I86").Value = 2841.5
$ws.Range("J86").Value = 2841.5715
$ws.Range("K86").Value = 2841.5
$ws.Range("L86").Value = 2841.5715
$ws.Range("M86").Value = -1718.5
$ws.Range("N86").Value = -5087.5715
$ws.Range("H89").Value = 2841.5173
$ws.Range("I89").Value = 2841.5
$ws.Range("J89").Value = 2841.5715
$ws.Range("K89").Value = 14207.5
$ws.Range("L89").Value = 14207.8575
$ws.Range("M89").Value = -8591.5
$ws.Range("N89").Value = -25439.8575
$ws.Range("H132").Value = 90998.39999999999
$ws.Range("J132").Value = 90998.39999999999
$ws.Range("L132").Value = 90998.39999999999
$ws.Range("N132").Value = -101118.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2846846.2
$ws.Range("I31").Value = 1880.8
$ws.Range("J31").Value = 4318380
$ws.Range("K31").Value = 1880.8
$ws.Range("L31").Value = 4318380
$ws.Range("M31").Value = -1585.8
$ws.Range("N31").Value = -4318970
$ws.Range("H34").Value = 2846846.2
$ws.Range("I34").Value = 1880.8
$ws.Range("J34").Value = 4318380
$ws.Range("K34").Value = 1880.8
$ws.Range("L34").Value = 4318380
$ws.Range("M34").Value = -1678.8
$ws.Range("N34").Value = -4318784
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H82").Value = 15001
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 15001
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H107").Value = 924.5
$ws.Range("I107").Value = 732.2222
$ws.Range("K107").Value = 732.2222
$ws.Range("M107").Value = 1187.7778
$ws.Range("H132").Value = 3656.6667
$ws.Range("I132").Value = 2886.6428
$ws.Range("J132").Value = 6351.75
$ws.Range("K132").Value = 8659.928400000001
$ws.Range("L132").Value = 19055.25
$ws.Range("M132").Value = -6129.928400000001
$ws.Range("N132").Value = -24115.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 87.19
$ws.Range("I4").Value = 68.125
$ws.Range("K4").Value = 204.375
$ws.Range("M4").Value = -92.375
$ws.Range("H80").Value = 4753
$ws.Range("J80").Value = 4753
$ws.Range("L80").Value = 14259
$ws.Range("N80").Value = -16131
$ws.Range("H83").Value = 4753
$ws.Range("J83").Value = 4753
$ws.Range("L83").Value = 42777
$ws.Range("N83").Value = -52137
$ws.Range("H92").Value = 1236.5
$ws.Range("I92").Value = 1998
$ws.Range("J92").Value = 475
$ws.Range("K92").Value = 5994
$ws.Range("L92").Value = 1425
$ws.Range("M92").Value = -4746
$ws.Range("N92").Value = -3921
$ws.Range("H103").Value = 770
$ws.Range("J103").Value = 560.8333
$ws.Range("L103").Value = 1682.4999
$ws.Range("N103").Value = -3440.4999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H126").Value = 31255428
$ws.Range("I126").Value = 50002280
$ws.Range("J126").Value = 10674.333
$ws.Range("K126").Value = 150006840
$ws.Range("L126").Value = 32022.999
$ws.Range("M126").Value = -150004370
$ws.Range("N126").Value = -36962.999
$ws.Range("H132").Value = 1664.079
$ws.Range("I132").Value = 1584.3055
$ws.Range("K132").Value = 4752.916499999999
$ws.Range("M132").Value = -2222.916499999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 37136.75
$ws.Range("J134").Value = 37136.75
$ws.Range("L134").Value = 111410.25
$ws.Range("N134").Value = -116480.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2363.4443
$ws.Range("I7").Value = 2181.8572
$ws.Range("K7").Value = 2181.8572
$ws.Range("M7").Value = -2069.8572
$ws.Range("H18").Value = 25000
$ws.Range("J18").Value = 25000
$ws.Range("L18").Value = 25000
$ws.Range("N18").Value = -25344
$ws.Range("H68").Value = 3011
$ws.Range("I68").Value = 3111.2
$ws.Range("J68").Value = 2844
$ws.Range("K68").Value = 3111.2
$ws.Range("L68").Value = 2844
$ws.Range("M68").Value = -2362.2
$ws.Range("N68").Value = -4342
$ws.Range("H71").Value = 3011
$ws.Range("I71").Value = 3111.2
$ws.Range("J71").Value = 2844
$ws.Range("K71").Value = 15556
$ws.Range("L71").Value = 14220
$ws.Range("M71").Value = -11812
$ws.Range("N71").Value = -21708
$ws.Range("H122").Value = 3479
$ws.Range("I122").Value = 3479
$ws.Range("K122").Value = 10437
$ws.Range("M122").Value = -7987
$ws.Range("H126").Value = 2363.4443
$ws.Range("I126").Value = 2181.8572
$ws.Range("K126").Value = 6545.571599999999
$ws.Range("M126").Value = -4075.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 3337333.2
$ws.Range("I40").Value = 5000000
$ws.Range("J40").Value = 2506000
$ws.Range("K40").Value = 5000000
$ws.Range("L40").Value = 2506000
$ws.Range("M40").Value = -4999851
$ws.Range("N40").Value = -2506298
$ws.Range("H82").Value = 100301
$ws.Range("J82").Value = 100301
$ws.Range("L82").Value = 100301
$ws.Range("N82").Value = -101067
$ws.Range("H85").Value = 100301
$ws.Range("J85").Value = 100301
$ws.Range("L85").Value = 100301
$ws.Range("N85").Value = -102953
$ws.Range("H136").Value = 8774457
$ws.Range("I136").Value = 10103557
$ws.Range("J136").Value = 2401.2
$ws.Range("K136").Value = 30310671
$ws.Range("L136").Value = 7203.599999999999
$ws.Range("M136").Value = -30308121
$ws.Range("N136").Value = -12303.6
